$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 836.53845
$ws.Range("I5").Value = 973
$ws.Range("J5").Value = 86
$ws.Range("K5").Value = 973
$ws.Range("L5").Value = 86
$ws.Range("M5").Value = -858
$ws.Range("N5").Value = -316
$ws.Range("H28").Value = 547.5454999999999
$ws.Range("I28").Value = 547.5454999999999
$ws.Range("K28").Value = 547.5454999999999
$ws.Range("M28").Value = -62.54549999999995
$ws.Range("H55").Value = 325.11765
$ws.Range("I55").Value = 309.07144
$ws.Range("J55").Value = 400
$ws.Range("K55").Value = 309.07144
$ws.Range("L55").Value = 400
$ws.Range("M55").Value = -95.07144
$ws.Range("N55").Value = -828
$ws.Range("H64").Value = 3567.6667
$ws.Range("I64").Value = 3274.4695
$ws.Range("J64").Value = 4593.857
$ws.Range("K64").Value = 3274.4695
$ws.Range("L64").Value = 4593.857
$ws.Range("M64").Value = -3026.4695
$ws.Range("N64").Value = -5089.857
$ws.Range("H67").Value = 3567.6667
$ws.Range("I67").Value = 3274.4695
$ws.Range("J67").Value = 4593.857
$ws.Range("K67").Value = 3274.4695
$ws.Range("L67").Value = 4593.857
$ws.Range("M67").Value = -2416.4695
$ws.Range("N67").Value = -6309.857
$ws.Range("H116").Value = 84452.69500000001
$ws.Range("I116").Value = 108490.5
$ws.Range("J116").Value = 4326.6665
$ws.Range("K116").Value = 108490.5
$ws.Range("L116").Value = 4326.6665
$ws.Range("M116").Value = -105048.5
$ws.Range("N116").Value = -11210.6665
$ws.Range("H132").Value = 2992.5
$ws.Range("I132").Value = 1355.4445
$ws.Range("J132").Value = 11177.777
$ws.Range("K132").Value = 4066.3335
$ws.Range("L132").Value = 33533.331
$ws.Range("M132").Value = -1536.3335
$ws.Range("N132").Value = -38593.331
$ws.Range("H137").Value = 5388.7715
$ws.Range("I137").Value = 5911.9653
$ws.Range("J137").Value = 2860
$ws.Range("K137").Value = 17735.8959
$ws.Range("L137").Value = 8580
$ws.Range("M137").Value = -15185.8959
$ws.Range("N137").Value = -13680
$ws.Range("H138").Value = 1892.3086
$ws.Range("I138").Value = 748.1579
$ws.Range("J138").Value = 2903.4187
$ws.Range("K138").Value = 2244.4737
$ws.Range("L138").Value = 8710.256100000001
$ws.Range("M138").Value = 2895.5263
$ws.Range("N138").Value = -18990.2561
$ws.Range("H141").Value = 2598.25
$ws.Range("I141").Value = 2213.3333
$ws.Range("J141").Value = 3753
$ws.Range("K141").Value = 6639.999899999999
$ws.Range("L141").Value = 11259
$ws.Range("M141").Value = -1459.999899999999
$ws.Range("N141").Value = -21619

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4003.5518
$ws.Range("I61").Value = 4100.0835
$ws.Range("J61").Value = 3935.4119
$ws.Range("K61").Value = 4100.0835
$ws.Range("L61").Value = 3935.4119
$ws.Range("M61").Value = -3888.0835
$ws.Range("N61").Value = -4359.4119
$ws.Range("H63").Value = 3025
$ws.Range("I63").Value = 3100
$ws.Range("J63").Value = 2875
$ws.Range("K63").Value = 3100
$ws.Range("L63").Value = 2875
$ws.Range("M63").Value = -2414
$ws.Range("N63").Value = -4247
$ws.Range("H66").Value = 3025
$ws.Range("I66").Value = 3100
$ws.Range("J66").Value = 2875
$ws.Range("K66").Value = 15500
$ws.Range("L66").Value = 14375
$ws.Range("M66").Value = -12068
$ws.Range("N66").Value = -21239
$ws.Range("H74").Value = 2094.6099
$ws.Range("I74").Value = 1301.2069
$ws.Range("J74").Value = 4012
$ws.Range("K74").Value = 1301.2069
$ws.Range("L74").Value = 4012
$ws.Range("M74").Value = -427.2068999999999
$ws.Range("N74").Value = -5760
$ws.Range("H77").Value = 2094.6099
$ws.Range("I77").Value = 1301.2069
$ws.Range("J77").Value = 4012
$ws.Range("K77").Value = 6506.0345
$ws.Range("L77").Value = 20060
$ws.Range("M77").Value = -2138.0345
$ws.Range("N77").Value = -28796
$ws.Range("H122").Value = 2373.5527
$ws.Range("I122").Value = 3614.7646
$ws.Range("J122").Value = 1368.762
$ws.Range("K122").Value = 10844.2938
$ws.Range("L122").Value = 4106.286
$ws.Range("M122").Value = -8394.293799999999
$ws.Range("N122").Value = -9006.286
$ws.Range("H134").Value = 34561.727
$ws.Range("J134").Value = 34561.727
$ws.Range("L134").Value = 34561.727
$ws.Range("N134").Value = -44701.727
$ws.Range("H136").Value = 4003.5518
$ws.Range("I136").Value = 4100.0835
$ws.Range("J136").Value = 3935.4119
$ws.Range("K136").Value = 12300.2505
$ws.Range("L136").Value = 11806.2357
$ws.Range("M136").Value = -9750.250499999998
$ws.Range("N136").Value = -16906.2357

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 362
$ws.Range("I22").Value = 322.4
$ws.Range("J22").Value = 560
$ws.Range("K22").Value = 322.4
$ws.Range("L22").Value = 560
$ws.Range("M22").Value = -149.4
$ws.Range("N22").Value = -906
$ws.Range("H99").Value = 8657755
$ws.Range("I99").Value = 2963082.8
$ws.Range("J99").Value = 33334666
$ws.Range("K99").Value = 2963082.8
$ws.Range("L99").Value = 33334666
$ws.Range("M99").Value = -2961584.8
$ws.Range("N99").Value = -33337662
$ws.Range("H105").Value = 5557514
$ws.Range("I105").Value = 2017.8572
$ws.Range("J105").Value = 25001750
$ws.Range("K105").Value = 2017.8572
$ws.Range("L105").Value = 25001750
$ws.Range("M105").Value = -270.8571999999999
$ws.Range("N105").Value = -25005244

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2273.3674
$ws.Range("I31").Value = 1437.0541
$ws.Range("J31").Value = 4852
$ws.Range("K31").Value = 1437.0541
$ws.Range("L31").Value = 4852
$ws.Range("M31").Value = -1142.0541
$ws.Range("N31").Value = -5442
$ws.Range("H34").Value = 2273.3674
$ws.Range("I34").Value = 1437.0541
$ws.Range("J34").Value = 4852
$ws.Range("K34").Value = 1437.0541
$ws.Range("L34").Value = 4852
$ws.Range("M34").Value = -1235.0541
$ws.Range("N34").Value = -5256
$ws.Range("H86").Value = 4439.077
$ws.Range("J86").Value = 5634.222
$ws.Range("L86").Value = 5634.222
$ws.Range("N86").Value = -7880.222
$ws.Range("H89").Value = 4439.077
$ws.Range("J89").Value = 5634.222
$ws.Range("L89").Value = 28171.11
$ws.Range("N89").Value = -39403.11
$ws.Range("H105").Value = 913.0484
$ws.Range("I105").Value = 827.72
$ws.Range("J105").Value = 1268.5834
$ws.Range("K105").Value = 827.72
$ws.Range("L105").Value = 1268.5834
$ws.Range("M105").Value = 919.28
$ws.Range("N105").Value = -4762.5834
$ws.Range("H122").Value = 2218.4119
$ws.Range("I122").Value = 2590.818
$ws.Range("J122").Value = 1535.6666
$ws.Range("K122").Value = 7772.454000000001
$ws.Range("L122").Value = 4606.9998
$ws.Range("M122").Value = -5322.454000000001
$ws.Range("N122").Value = -9506.9998
$ws.Range("H134").Value = 1523.0646
$ws.Range("I134").Value = 808.65216
$ws.Range("K134").Value = 2425.95648
$ws.Range("M134").Value = 109.0435200000002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3795.1853
$ws.Range("I102").Value = 2713.4736
$ws.Range("J102").Value = 6364.25
$ws.Range("K102").Value = 2713.4736
$ws.Range("L102").Value = 6364.25
$ws.Range("M102").Value = -1091.4736
$ws.Range("N102").Value = -9608.25
$ws.Range("H105").Value = 35671
$ws.Range("J105").Value = 35671
$ws.Range("L105").Value = 35671
$ws.Range("N105").Value = -42659
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H126").Value = 1820.96
$ws.Range("I126").Value = 1963.3334
$ws.Range("J126").Value = 1689.5385
$ws.Range("K126").Value = 5890.0002
$ws.Range("L126").Value = 5068.6155
$ws.Range("M126").Value = -3420.0002
$ws.Range("N126").Value = -10008.6155

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3585.3845
$ws.Range("I7").Value = 4275
$ws.Range("J7").Value = 2482
$ws.Range("K7").Value = 4275
$ws.Range("L7").Value = 2482
$ws.Range("M7").Value = -4163
$ws.Range("N7").Value = -2706
$ws.Range("H106").Value = 28265.25
$ws.Range("J106").Value = 28265.25
$ws.Range("L106").Value = 28265.25
$ws.Range("N106").Value = -30789.25
$ws.Range("H126").Value = 3585.3845
$ws.Range("I126").Value = 4275
$ws.Range("J126").Value = 2482
$ws.Range("K126").Value = 12825
$ws.Range("L126").Value = 7446
$ws.Range("M126").Value = -10355
$ws.Range("N126").Value = -12386
$ws.Range("H132").Value = 10644876
$ws.Range("I132").Value = 20002376
$ws.Range("J132").Value = 11354.454
$ws.Range("K132").Value = 60007128
$ws.Range("L132").Value = 34063.362
$ws.Range("M132").Value = -60004598
$ws.Range("N132").Value = -39123.362

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 35196
$ws.Range("J105").Value = 35196
$ws.Range("L105").Value = 35196
$ws.Range("N105").Value = -42184
$ws.Range("H122").Value = 1961
$ws.Range("I122").Value = 1873.7858
$ws.Range("J122").Value = 2368
$ws.Range("K122").Value = 5621.357400000001
$ws.Range("L122").Value = 7104
$ws.Range("M122").Value = -3171.357400000001
$ws.Range("N122").Value = -12004
$ws.Range("H132").Value = 18520246
$ws.Range("I132").Value = 25641984
$ws.Range("K132").Value = 76925952
$ws.Range("M132").Value = -76923422
$ws.Range("H135").Value = 38422.855
$ws.Range("J135").Value = 38422.855
$ws.Range("L135").Value = 38422.855
$ws.Range("N135").Value = -48562.855
